$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# Set the Experimental value (row 7, column B) to the literal text "true".
# A direct Value assignment of "true" gets auto-coerced to a Boolean by the
# engine (mirrors Excel's literal-boolean recognition), so we build the text
# via a helper formula and paste the computed result back as a value - this
# keeps the cell's text type (rather than boolean) and preserves the
# existing cell style.
$helper = $ws.Range("Z1")
$helper.Formula = "=""true"""
$helper.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
